$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet/tab
$ws.Name = "위키 데이터베이스 기능 목록표"

# Complete row 6 (getDocumentRelDocIndex)
$ws.Range("C6").Value = "JSON reldocIndex"
$ws.Range("D6").Value = "[관련있는,문서]"
$ws.Range("E6").Value = "docID 문서의 상위/관련 문서 리스트를 가져 옵니다."

# Complete row 7 (getDocumentFrameList)
$ws.Range("B7").Value = "docID"
$ws.Range("C7").Value = "JSON frameList"
$ws.Range("D7").Value = "[Frame:임시,Frame:틀1]"
$ws.Range("E7").Value = "docID 문서의 틀 리스트를 가져옵니다."

# Row 8 - getDocumentFullData
$ws.Range("A8").Value = "getDocumentFullData"
$ws.Range("B8").Value = "docID"
$ws.Range("C8").Value = "JSON Parsed Document Data"
$ws.Range("D8").Value = "타 시트 참조"
$ws.Range("E8").Value = "docID 문서의 파싱된 데이터를 가져 옵니다."

# Row 9 - getDocumentRawData
$ws.Range("A9").Value = "getDocumentRawData"
$ws.Range("B9").Value = "docID"
$ws.Range("C9").Value = "JSON Raw(Non-parsed) Document Data"
$ws.Range("D9").Value = "타 시트 참조"
$ws.Range("E9").Value = "docID 문서의 데이터를 가져옵니다."

# Row 10 - isDocumentExists
$ws.Range("A10").Value = "isDocumentExists"
$ws.Range("B10").Value = "docID"
$ws.Range("C10").Value = "Boolean docExists"
$ws.Range("D10").Value = $true
$ws.Range("E10").Value = "문서가 존재하는지/아닌지를 반환합니다."

# Row 11 - getFrameData (A11 uses a distinct font, family=3)
$ws.Range("A11").Value = "getFrameData"
$ws.Range("A11").Font.Family = 3
$ws.Range("B11").Value = "frameID, frameTier"
$ws.Range("C11").Value = "String Data"
$ws.Range("D11").Value = "`"'''이 틀은 임시 틀입니다!'''`""
$ws.Range("E11").Value = "FrameTier등급의 FrameID 틀을 가져옵니다."

# Row 12 - testFunction
$ws.Range("A12").Value = "testFunction"
$ws.Range("B12").Value = "없음"
$ws.Range("C12").Value = "없음"
$ws.Range("D12").Value = "없음"
$ws.Range("E12").Value = "테스트 값을 입력합니다."

# Page setup: paper size 9 (A4), portrait orientation
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Update selection to match final state
$ws.Range("B15").Select() | Out-Null
